{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// 1) Lengthen the two blank-line underscores in the \"Copies served by Dep. Clerk\" line.\n// 2) Tighten one space in \"Prosecutor's Office\" line and add a trailing semicolon.\n// 3) Insert a brand-new paragraph (same tab stops / font / size) right after the\n//    \"Prosecutor's Office ...\" paragraph with the Community Control / County Jail COS line.\n\nconst OLD_COPIES = \"Copies served by Dep. Clerk ___________ on the following date ___________ to:\";\nconst NEW_COPIES = \"Copies served by Dep. Clerk ___________________________ on the following date ____________________ to:\";\n\nconst OLD_COS = \"Prosecutor\u2019s Office: PS     OM     EM; Defendant\u2019s Attorney: PS     OM     EM; {{ defendant.first_name }} {{ defendant.last_name}}: PS     OM     EM\";\nconst NEW_COS = \"Prosecutor\u2019s Office: PS    OM     EM; Defendant\u2019s Attorney: PS     OM     EM; {{ defendant.first_name }} {{ defendant.last_name}}: PS     OM     EM;\";\n\nconst NEW_PARA_TEXT = \"{% if community_control.ordered is true or bond_conditions.monitoring is true %}Community Control: PS    EM;{% endif %}{% if jail_terms.ordered is true or apply_jtc == \u2018Sentence\u2019 %}County Jail: PS   EM;{% endif %}\";\n\n// --- Edit 1: extend the underscores on the \"Copies served by Dep. Clerk\" line ---\nconst copiesResults = context.document.body.search(OLD_COPIES, { matchCase: true });\ncopiesResults.load(\"items\");\nawait context.sync();\n\nif (copiesResults.items.length === 0) {\n  throw new Error(\"Could not find the 'Copies served by Dep. Clerk' paragraph text.\");\n}\ncopiesResults.items[0].insertText(NEW_COPIES, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 2: tighten spacing + trailing semicolon on the \"Prosecutor's Office\" line ---\nconst cosResults = context.document.body.search(OLD_COS, { matchCase: true });\ncosResults.load(\"items\");\nawait context.sync();\n\nif (cosResults.items.length === 0) {\n  throw new Error(\"Could not find the 'Prosecutor's Office ...' paragraph text.\");\n}\nconst cosRange = cosResults.items[0];\ncosRange.insertText(NEW_COS, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 3: insert a new paragraph right after the \"Prosecutor's Office\" paragraph ---\n// Re-search on the (now updated) text so we anchor on the paragraph itself.\nconst cosResults2 = context.document.body.search(NEW_COS, { matchCase: true });\ncosResults2.load(\"items\");\nawait context.sync();\n\nif (cosResults2.items.length === 0) {\n  throw new Error(\"Could not re-find the updated 'Prosecutor's Office ...' text.\");\n}\nconst cosParagraph = cosResults2.items[0].paragraphs.getFirst();\nconst newParagraph = cosParagraph.insertParagraph(NEW_PARA_TEXT, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument / $d is the live document.\n#\n# 1) Lengthen the two blank-line underscores in the \"Copies served by Dep. Clerk\" line.\n# 2) Tighten one space in \"Prosecutor's Office\" line and add a trailing semicolon.\n# 3) Insert a brand-new paragraph (same tab stops / font / size) right after the\n#    \"Prosecutor's Office ...\" paragraph with the Community Control / County Jail COS line.\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: extend the underscores on the \"Copies served by Dep. Clerk\" line ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = 'Copies served by Dep. Clerk ___________ on the following date ___________ to:'\n$find1.Replacement.Text = 'Copies served by Dep. Clerk ___________________________ on the following date ____________________ to:'\n$found1 = $find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find1.Replacement.Text, $wdReplaceOne)\nif (-not $found1) {\n    throw \"Could not find/replace the 'Copies served by Dep. Clerk' paragraph text.\"\n}\n\n# --- Edit 2: tighten spacing + trailing semicolon on the \"Prosecutor's Office\" line ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = 'Prosecutor\u2019s Office: PS     OM     EM; Defendant\u2019s Attorney: PS     OM     EM; {{ defendant.first_name }} {{ defendant.last_name}}: PS     OM     EM'\n$find2.Replacement.Text = 'Prosecutor\u2019s Office: PS    OM     EM; Defendant\u2019s Attorney: PS     OM     EM; {{ defendant.first_name }} {{ defendant.last_name}}: PS     OM     EM;'\n$found2 = $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceOne)\nif (-not $found2) {\n    throw \"Could not find/replace the 'Prosecutor's Office ...' paragraph text.\"\n}\n\n# --- Edit 3: insert a new paragraph right after the \"Prosecutor's Office\" paragraph ---\n# Re-locate the paragraph that now holds the updated COS text so the new\n# paragraph lands in the right place and inherits its tab stops / font / size.\n$targetText = 'Prosecutor\u2019s Office: PS    OM     EM; Defendant\u2019s Attorney: PS     OM     EM; {{ defendant.first_name }} {{ defendant.last_name}}: PS     OM     EM;'\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText) {\n        $targetParagraph = $p\n        break\n    }\n}\nif ($null -eq $targetParagraph) {\n    throw \"Could not re-find the updated 'Prosecutor's Office ...' paragraph.\"\n}\n\n$targetParagraph.Range.InsertParagraphAfter()\n$newParagraph = $targetParagraph.Next()\n$newParagraph.Range.Text = '{% if community_control.ordered is true or bond_conditions.monitoring is true %}Community Control: PS    EM;{% endif %}{% if jail_terms.ordered is true or apply_jtc == \u2018Sentence\u2019 %}County Jail: PS   EM;{% endif %}'\n"}
